$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 currently holds "Password"; change it to lowercase "password".
# B3 (which also reads "Password") keeps its own text unchanged - the engine
# tracks each cell's string value independently, so it is unaffected by this.
$ws.Range("B1").Value = "password"

# Update the selected cell to match the saved view state in the target file.
$ws.Range("C11").Select()
